$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 213
$ws.Range("C1").Value = 123.213122
$ws.Range("D1").Value = 52
$ws.Range("E1").Value = 231
$ws.Range("F1").Value = 231
$ws.Range("G1").Value = "Hello world"

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 2

$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 2
$ws.Range("C7").Value = 5

[void]$ws.Range("A1").Select()
